$wb = $excel.ActiveWorkbook

# The handoff for the first tracked file failed, so the report status
# changes everywhere "Not yet handed off" was shown.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handoff failed"
$overview.Range("C2").Value = "Handoff failed"

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remove just the hyperlink anchored at C2 (handoff failed -> no handoff file link)
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq '$C$2') {
            $hl.Delete()
        }
    }

    # Row 2 reflects a failed handoff attempt:
    #  - Status -> "Handoff failed"
    #  - Latest Handoff File (C2) cleared entirely (no file, no hyperlink)
    #  - Latest Handoff Datetime (D2) reset to the zero-date
    #  - Handoff Reason (H2) -> "Ignored"
    $ws.Range("B2").Value = "Handoff failed"
    $ws.Range("C2").Clear()
    $ws.Range("D2").Value = "0001-01-01 00:00:00"
    $ws.Range("H2").Value = "Ignored"
}
